# TRIMAZKON_address_list.xlsx - "ip_address_list" sheet rewrite
#
# The sheet holds a flat address-book table (columns A..E: Name, IP,
# Netmask, Notes, Fav-flag) as plain text cells. This commit rewrites
# most of the rows: several entries are renamed/re-IP'd, some long
# "Notes" blobs get small in-place edits, a couple of duplicate rows are
# appended, and the table grows from 10 to 13 rows.
#
# Simplest faithful way to reproduce that with the Excel object model is
# to just (re)write every cell A1:E13 to its final value - this both
# updates existing rows and creates the new trailing ones in one pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ip_address_list")

# Columns A and E sometimes hold purely-numeric-looking text ("47", "0",
# "1", ...). Excel's normal .Value assignment would auto-coerce those to
# real numbers, but the source file stores them as text (t="inlineStr").
# Pre-formatting the range as text keeps every assignment below a string.
$ws.Range("A1:E13").NumberFormat = "@"

# Row 1
$ws.Range("A1").Value = "474 B_Austin"
$ws.Range("B1").Value = "10.96.205.175"
$ws.Range("C1").Value = "255.255.255.0"
$ws.Range("D1").Value = "PC:`t10.96.aoj"
$ws.Range("E1").Value = "1"

# Row 2
$ws.Range("A2").Value = "529_Witte55"
$ws.Range("B2").Value = "192.168.0.240"
$ws.Range("C2").Value = "255.255.255.0"
$ws.Range("D2").Value = "P"
$ws.Range("E2").Value = "0"

# Row 3 (unchanged by the commit, rewritten here for completeness)
$ws.Range("A3").Value = "474 B_Austin (1)"
$ws.Range("B3").Value = "10.96.205.175"
$ws.Range("C3").Value = "255.255.255.0"
$ws.Range("D3").Value = "PC:`t10.96.205.175`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.20`n-----------------------------------------`nuser:JHV_Vision, omron `nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"
$ws.Range("E3").Value = "1"

# Row 4
$ws.Range("A4").Value = "474 B_Austin (2)"
$ws.Range("B4").Value = "10.96.205.175"
$ws.Range("C4").Value = "255.255.255.0"
$ws.Range("D4").Value = "10.96.205.1`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.20`n-----------------------------------------`nuser:JHV_Vision, omron llllllllllllll`nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"
$ws.Range("E4").Value = "1"

# Row 5 (new)
$ws.Range("A5").Value = "474 B_Austin (2) (1)"
$ws.Range("B5").Value = "10.96.205.175"
$ws.Range("C5").Value = "255.255.255.0"
$ws.Range("D5").Value = "10.96.205.1`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.20`n-----------------------------------------`nuser:JHV_Vision, omron llllllllllllll`nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"
$ws.Range("E5").Value = "1"

# Row 6 (new)
$ws.Range("A6").Value = "529_Witte"
$ws.Range("B6").Value = "192.168.0.240"
$ws.Range("C6").Value = "255.255.255.0"
$ws.Range("D6").Value = "PC"
$ws.Range("E6").Value = "0"

# Row 7 (new)
$ws.Range("A7").Value = "47"
$ws.Range("B7").Value = "10.96.205.175"
$ws.Range("C7").Value = "255.255.255.0"
$ws.Range("D7").Value = "PC:`t10.96.205."
$ws.Range("E7").Value = "1"

# Row 8 (new)
$ws.Range("A8").Value = "518_Val"
$ws.Range("B8").Value = "192.168.208.242"
$ws.Range("C8").Value = "255.255.255.0"
$ws.Range("D8").ClearContents()
$ws.Range("E8").Value = "0"

# Row 9 (new)
$ws.Range("A9").Value = "518_Valeo II"
$ws.Range("B9").Value = "192.168.1.243"
$ws.Range("C9").Value = "255.255.255.0"
$ws.Range("D9").Value = "ssssssss"
$ws.Range("E9").Value = "0"

# Row 10
$ws.Range("A10").Value = "527_"
$ws.Range("B10").Value = "10.101.28.176"
$ws.Range("C10").Value = "255.255.255.0"
$ws.Range("D10").Value = "PC:`t10.96.20"
$ws.Range("E10").Value = "0"

# Row 11 (new)
$ws.Range("A11").Value = "Dom"
$ws.Range("B11").Value = "192.168.1.131"
$ws.Range("C11").Value = "255.255.255.0"
$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = "1"

# Row 12 (new)
$ws.Range("A12").Value = "474 B_A"
$ws.Range("B12").Value = "10.96.205.175"
$ws.Range("C12").Value = "255.255.255.0"
$ws.Range("D12").Value = "dfddddddddddddddddd`nadf`nafd`nafsdfaadfs"
$ws.Range("E12").Value = "0"

# Row 13 (new)
$ws.Range("A13").Value = "Domac"
$ws.Range("B13").Value = "192.168.1.13"
$ws.Range("C13").Value = "255.255.255.0"
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = "0"
